# Activity name changes for samsung
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): add new "Notification" column header in G1 ---
$ws.Cells.Item(1, 7).Value = "Notification"

# --- Row 2: CallHomeforLess OPT IN ---
$ws.Cells.Item(2, 2).Value = "CallHomeforLess - OPT IN"
$ws.Cells.Item(2, 3).Value = "CallHomeforLess"
$ws.Cells.Item(2, 4).Value = "*135*30#"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = "device1"
$ws.Cells.Item(2, 7).Value = "The call home for less offer is active! Call at only 0.6 fils per sec, with a setup fee of AED 1 per call for calls to India, Bangladesh, Pakistan, Afghanistan (except numbers starting with 009378), Egypt, China, Iran, Nepal, Nigeria & Indonesia."

# --- Row 3: CallHomeforLess OPT OUT ---
$ws.Cells.Item(3, 2).Value = "CallHomeforLess - OPT OUT"
$ws.Cells.Item(3, 3).Value = "CallHomeforLess"
$ws.Cells.Item(3, 4).Value = "*135*30#"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = "device1"
$ws.Cells.Item(3, 7).Value = "As requested, you've been opted out from the New Call Home For Less promotion."

# --- Formatting for the new Notification column (G2:G3) ---
$notifRange = $ws.Range("G2:G3")
$notifRange.Font.Name = "Ericsson Hilda"
$notifRange.Font.Size = 9
$notifRange.Borders.LineStyle = 1
$notifRange.HorizontalAlignment = -4131
$notifRange.VerticalAlignment = -4160

# G3 additionally wraps text
$ws.Cells.Item(3, 7).WrapText = $true

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 20.26
$ws.Columns.Item(7).ColumnWidth = 175.53

# --- Selection, matching the saved view state ---
$ws.Range("C1").Select()
